$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Octubre de 2020 a las 18:00"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 8165007
$ws.Cells.Item(4, 3).Value = 14964
$ws.Cells.Item(4, 4).Value = 5290510
$ws.Cells.Item(4, 5).Value = 2652405
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 249
$ws.Cells.Item(4, 8).Value = 222092

# Row 5
$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 7349290
$ws.Cells.Item(5, 3).Value = 44220
$ws.Cells.Item(5, 4).Value = 6425716
$ws.Cells.Item(5, 5).Value = 811848
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 415
$ws.Cells.Item(5, 8).Value = 111726

# Row 6
$ws.Cells.Item(6, 1).Value = "Brasil"
$ws.Cells.Item(6, 2).Value = 5142003
$ws.Cells.Item(6, 3).Value = 505
$ws.Cells.Item(6, 4).Value = 4568813
$ws.Cells.Item(6, 5).Value = 421410
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 151780

# Row 15
$ws.Cells.Item(15, 1).Value = "Reino Unido"
$ws.Cells.Item(15, 2).Value = 673622
$ws.Cells.Item(15, 3).Value = 18980
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 138
$ws.Cells.Item(15, 8).Value = 43293

# Row 17
$ws.Cells.Item(17, 1).Value = "Chile"
$ws.Cells.Item(17, 2).Value = 486496
$ws.Cells.Item(17, 3).Value = 1124
$ws.Cells.Item(17, 4).Value = 459536
$ws.Cells.Item(17, 5).Value = 13526
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 19
$ws.Cells.Item(17, 8).Value = 13434

# Row 20
$ws.Cells.Item(20, 1).Value = "Italia"
$ws.Cells.Item(20, 2).Value = 381602
$ws.Cells.Item(20, 3).Value = 8804
$ws.Cells.Item(20, 4).Value = 245964
$ws.Cells.Item(20, 5).Value = 99266
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 83
$ws.Cells.Item(20, 8).Value = 36372

# Row 30
$ws.Cells.Item(30, 1).Value = "Canada"
$ws.Cells.Item(30, 2).Value = 191137
$ws.Cells.Item(30, 3).Value = 1750
$ws.Cells.Item(30, 4).Value = 160130
$ws.Cells.Item(30, 5).Value = 21310
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 33
$ws.Cells.Item(30, 8).Value = 9697

# Row 48
$ws.Cells.Item(48, 1).Value = "Guatemala"
$ws.Cells.Item(48, 2).Value = 99765
$ws.Cells.Item(48, 3).Value = 671
$ws.Cells.Item(48, 4).Value = 88931
$ws.Cells.Item(48, 5).Value = 7381
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 23
$ws.Cells.Item(48, 8).Value = 3453

# Row 59
$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 65076
$ws.Cells.Item(59, 3).Value = 652
$ws.Cells.Item(59, 4).Value = 46010
$ws.Cells.Item(59, 5).Value = 17536
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 16
$ws.Cells.Item(59, 8).Value = 1530

# Row 64
$ws.Cells.Item(64, 1).Value = "Singapur"
$ws.Cells.Item(64, 2).Value = 57892
$ws.Cells.Item(64, 3).Value = 3
$ws.Cells.Item(64, 4).Value = 57764
$ws.Cells.Item(64, 5).Value = 100
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 28

# Row 80
$ws.Cells.Item(80, 1).Value = "Jordania"
$ws.Cells.Item(80, 2).Value = 33009
$ws.Cells.Item(80, 3).Value = 2459
$ws.Cells.Item(80, 4).Value = 6565
$ws.Cells.Item(80, 5).Value = 26162
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 25
$ws.Cells.Item(80, 8).Value = 282

# Row 81
$ws.Cells.Item(81, 1).Value = "Birmania"
$ws.Cells.Item(81, 2).Value = 32351
$ws.Cells.Item(81, 3).Value = 1026
$ws.Cells.Item(81, 4).Value = 14706
$ws.Cells.Item(81, 5).Value = 16880
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 33
$ws.Cells.Item(81, 8).Value = 765

# Row 82
$ws.Cells.Item(82, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(82, 2).Value = 32224
$ws.Cells.Item(82, 3).Value = 569
$ws.Cells.Item(82, 4).Value = 24356
$ws.Cells.Item(82, 5).Value = 6896
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 14
$ws.Cells.Item(82, 8).Value = 972

# Row 83
$ws.Cells.Item(83, 1).Value = "El Salvador"
$ws.Cells.Item(83, 2).Value = 31061
$ws.Cells.Item(83, 3).Value = 295
$ws.Cells.Item(83, 4).Value = 26311
$ws.Cells.Item(83, 5).Value = 3842
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 4
$ws.Cells.Item(83, 8).Value = 908

# Row 88
$ws.Cells.Item(88, 1).Value = "Grecia"
$ws.Cells.Item(88, 2).Value = 23947
$ws.Cells.Item(88, 3).Value = 452
$ws.Cells.Item(88, 4).Value = 9989
$ws.Cells.Item(88, 5).Value = 13476
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 13
$ws.Cells.Item(88, 8).Value = 482

# Row 95
$ws.Cells.Item(95, 1).Value = "Albania"
$ws.Cells.Item(95, 2).Value = 16212
$ws.Cells.Item(95, 3).Value = 257
$ws.Cells.Item(95, 4).Value = 9864
$ws.Cells.Item(95, 5).Value = 5909
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 5
$ws.Cells.Item(95, 8).Value = 439

# Row 96
$ws.Cells.Item(96, 1).Value = "Noruega"
$ws.Cells.Item(96, 2).Value = 16050
$ws.Cells.Item(96, 3).Value = 97
$ws.Cells.Item(96, 4).Value = 11863
$ws.Cells.Item(96, 5).Value = 3909
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 278

# Row 99
$ws.Cells.Item(99, 1).Value = "Montenegro"
$ws.Cells.Item(99, 2).Value = 14672
$ws.Cells.Item(99, 3).Value = 211
$ws.Cells.Item(99, 4).Value = 10355
$ws.Cells.Item(99, 5).Value = 4096
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 4
$ws.Cells.Item(99, 8).Value = 221

# Row 110
$ws.Cells.Item(110, 1).Value = "Luxemburgo"
$ws.Cells.Item(110, 2).Value = 10244
$ws.Cells.Item(110, 3).Value = 214
$ws.Cells.Item(110, 4).Value = 8384
$ws.Cells.Item(110, 5).Value = 1727
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 133

# Row 111
$ws.Cells.Item(111, 1).Value = "Guayana Francesa"
$ws.Cells.Item(111, 2).Value = 10233
$ws.Cells.Item(111, 3).Value = 31
$ws.Cells.Item(111, 4).Value = 9894
$ws.Cells.Item(111, 5).Value = 270
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 69

# Row 112
$ws.Cells.Item(112, 1).Value = "Uganda"
$ws.Cells.Item(112, 2).Value = 10117
$ws.Cells.Item(112, 3).Value = 48
$ws.Cells.Item(112, 4).Value = 6725
$ws.Cells.Item(112, 5).Value = 3296
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = 96

# Row 115
$ws.Cells.Item(115, 1).Value = "Jamaica"
$ws.Cells.Item(115, 2).Value = 8067
$ws.Cells.Item(115, 3).Value = 78
$ws.Cells.Item(115, 4).Value = 3481
$ws.Cells.Item(115, 5).Value = 4426
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 9
$ws.Cells.Item(115, 8).Value = 160

# Row 116
$ws.Cells.Item(116, 1).Value = "Zimbabue"
$ws.Cells.Item(116, 2).Value = 8055
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 7640
$ws.Cells.Item(116, 5).Value = 184
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 231

# Row 126
$ws.Cells.Item(126, 1).Value = "Bahamas"
$ws.Cells.Item(126, 2).Value = 5385
$ws.Cells.Item(126, 3).Value = 194
$ws.Cells.Item(126, 4).Value = 3178
$ws.Cells.Item(126, 5).Value = 2095
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 3
$ws.Cells.Item(126, 8).Value = 112

# Row 127
$ws.Cells.Item(127, 1).Value = "Nicaragua"
$ws.Cells.Item(127, 2).Value = 5353
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 4225
$ws.Cells.Item(127, 5).Value = 974
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 154

# Row 128
$ws.Cells.Item(128, 1).Value = "Sri Lanka"
$ws.Cells.Item(128, 2).Value = 5219
$ws.Cells.Item(128, 3).Value = 49
$ws.Cells.Item(128, 4).Value = 3380
$ws.Cells.Item(128, 5).Value = 1826
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 13

# Row 129
$ws.Cells.Item(129, 1).Value = "Hong Kong"
$ws.Cells.Item(129, 2).Value = 5214
$ws.Cells.Item(129, 3).Value = 12
$ws.Cells.Item(129, 4).Value = 4943
$ws.Cells.Item(129, 5).Value = 166
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 105

# Row 170
$ws.Cells.Item(170, 1).Value = "San Marino"
$ws.Cells.Item(170, 2).Value = 759
$ws.Cells.Item(170, 3).Value = 18
$ws.Cells.Item(170, 4).Value = 685
$ws.Cells.Item(170, 5).Value = 32
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 42
